# Report interactions User in Excel
# Adds two detail rows (interaction events) to the "informe interaccion stand"
# worksheet: a stand-element click/interaction log for "User1" (Secpho /
# Engineering / CEO), one row where a "Soft" element was touched and contact
# info was given (contactame@gmail.com), and one row where a "Hard" element
# was touched with no contact info given.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("informe interaccion stand")

$interactionDate = 44172.6333333333

# Row 2: User1 / Secpho / Engineering / CEO / Soft / contactame@gmail.com
$ws.Cells.Item(2, 4).Value = $interactionDate
$ws.Cells.Item(2, 4).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2, 5).Value = "User1"
$ws.Cells.Item(2, 6).Value = "Secpho"
$ws.Cells.Item(2, 7).Value = "Engineering"
$ws.Cells.Item(2, 8).Value = "CEO"
$ws.Cells.Item(2, 9).Value = "Soft"
$ws.Cells.Item(2, 10).Value = "contactame@gmail.com"

# Row 3: User1 / Secpho / Engineering / CEO / Hard / (no contact given)
$ws.Cells.Item(3, 4).Value = $interactionDate
$ws.Cells.Item(3, 4).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3, 5).Value = "User1"
$ws.Cells.Item(3, 6).Value = "Secpho"
$ws.Cells.Item(3, 7).Value = "Engineering"
$ws.Cells.Item(3, 8).Value = "CEO"
$ws.Cells.Item(3, 9).Value = "Hard"
$ws.Cells.Item(3, 10).Value = ""
